$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("car inventory")

# Add "Miles / Year" formula in column I for row 2 first (its own, non-shared formula)
$ws.Range("I2").Formula = "=H2/G2"

# Then fill in the rest of the column (rows 3-53) as a separate shared formula group
$ws.Range("I3:I53").Formula = "=H3/G3"

# Update the active selection to J2, matching the saved workbook state
$ws.Range("J2").Select()
